# The deck ships two embedded themes:
#   ppt/theme/theme1.xml -> bound to the Slide Master, clrScheme "Integral"
#   ppt/theme/theme2.xml -> bound to the Notes Master,  clrScheme "Office"
# The target revision swaps the two themes' contents (their color schemes -
# font scheme and format scheme are already byte-identical between the two
# parts, so a full "Integral" <-> "Office" swap is really just a 12-slot
# theme-color swap). We re-point the Slide Master's theme color scheme from
# the "Integral" palette to the "Office" palette via the PowerPoint object
# model (Theme.ThemeColorScheme), which is the supported automation surface
# for rewriting a design's clrScheme entries in place.

function Convert-HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette ("Office Theme") in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

$slideMasterColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $slideMasterColors.Item($i).RGB = Convert-HexToRgbInt $officeTheme[$i - 1]
}

# Major/minor latin typefaces already match ("Arial") between the two
# themes, but set them explicitly so the design's font scheme stays in
# sync with the new color scheme.
$fontScheme = $p.SlideMaster.Theme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Arial"
$fontScheme.MinorFont.Latin = "Arial"
